$wb = $excel.ActiveWorkbook

# --- "workbooks" sheet: update the capture-test path string ---
$wsWorkbooks = $wb.Worksheets.Item("workbooks")
$wsWorkbooks.Range("A2").Value = "C:\Users\HP\git\Desktop-Framework\Finalized\2_Capture Test.xls"

# --- "Apps" sheet: update the notepad path entry (row for "notepad" app, B5) ---
$wsApps = $wb.Worksheets.Item("Apps")
$wsApps.Range("B5").Value = "C:\WINDOWS\system32\\notepad.exe"

# Update the selection shown in the "Apps" sheet view to B5
$wsApps.Range("B5").Select()

# Make "Apps" the active / selected sheet (tab) in the workbook
$wsApps.Activate()
